$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 681 (shifts the existing rows 681-711 down to 683-713,
# pulling the date-number-format style of column D down with them).
$ws.Rows.Item(681).Insert()
$ws.Rows.Item(681).Insert()

# New row 681: Acelga, "Primera" quality, week of 2023-08-09
$ws.Cells.Item(681, 1).Value = 8
$ws.Cells.Item(681, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(681, 3).Value = "Coquimbo"
$ws.Cells.Item(681, 4).Value = 45147
$ws.Cells.Item(681, 5).Value = 4
$ws.Cells.Item(681, 6).Value = 100112009
$ws.Cells.Item(681, 7).Value = "Acelga"
$ws.Cells.Item(681, 8).Value = "Sin especificar"
$ws.Cells.Item(681, 9).Value = "Primera"
$ws.Cells.Item(681, 10).Value = 2000
$ws.Cells.Item(681, 11).Value = 550
$ws.Cells.Item(681, 12).Value = 600
$ws.Cells.Item(681, 13).Value = 575
$ws.Cells.Item(681, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(681, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(681, 16).Value = 288
$ws.Cells.Item(681, 17).Value = 2
$ws.Cells.Item(681, 18).Value = "Hortaliza"

# New row 682: Acelga, "Segunda" quality, week of 2023-08-09
$ws.Cells.Item(682, 1).Value = 8
$ws.Cells.Item(682, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(682, 3).Value = "Coquimbo"
$ws.Cells.Item(682, 4).Value = 45147
$ws.Cells.Item(682, 5).Value = 4
$ws.Cells.Item(682, 6).Value = 100112009
$ws.Cells.Item(682, 7).Value = "Acelga"
$ws.Cells.Item(682, 8).Value = "Sin especificar"
$ws.Cells.Item(682, 9).Value = "Segunda"
$ws.Cells.Item(682, 10).Value = 1400
$ws.Cells.Item(682, 11).Value = 450
$ws.Cells.Item(682, 12).Value = 500
$ws.Cells.Item(682, 13).Value = 475
$ws.Cells.Item(682, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(682, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(682, 16).Value = 238
$ws.Cells.Item(682, 17).Value = 2
$ws.Cells.Item(682, 18).Value = "Hortaliza"
